$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = -0.4031925586089656
$ws.Range("J2").Value = 0.109469479763775
$ws.Range("K2").Value = -0.8233677731302166
$ws.Range("L2").Value = 2.776224066834431

$ws.Range("I12").Value = -0.2177642701949869
$ws.Range("J12").Value = 0.0311526845331437
$ws.Range("K12").Value = -0.5953334235215503
$ws.Range("L12").Value = 2.249093190457154

$ws.Range("I14").Value = -0.2966519979474677
$ws.Range("J14").Value = 0.07082151409614748
$ws.Range("K14").Value = -0.9037196749651489
$ws.Range("L14").Value = 2.727714024660366

$ws.Range("I18").Value = -0.5874877635030271
$ws.Range("J18").Value = 0.1412024805547104
$ws.Range("K18").Value = -0.0955827746949887
$ws.Range("L18").Value = 1.806117529550021
